# Commit: "Now parsing boolean signals as well!"
#
# 1) Rename the boolean Fault signal names in column B (rows 75-97, the
#    Protections_BMSFaults_Pack1/2/3 blocks) by appending a "_TF" suffix
#    so they're recognisable as true/false boolean signals.
# 2) Freeze the header row and move the active selection to E10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$faultRows = @(75, 76, 77, 78, 79, 80, 81, 83, 84, 85, 86, 87, 88, 89, 91, 92, 93, 94, 95, 96, 97)

foreach ($r in $faultRows) {
    $cell = $ws.Cells.Item($r, 2)
    $name = $cell.Value2
    if ($name -and -not $name.EndsWith("_TF")) {
        $cell.Value = $name + "_TF"
    }
}

# Freeze panes below row 1, and set the active selection to E10.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("E10").Select()
